$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Cases by Age Group"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = 268
$ws1.Range("B3").Value = 1273
$ws1.Range("B4").Value = 3465
$ws1.Range("B5").Value = 15149
$ws1.Range("B6").Value = 16759
$ws1.Range("B7").Value = 14665
$ws1.Range("B8").Value = 12304
$ws1.Range("B9").Value = 4430
$ws1.Range("B10").Value = 2972
$ws1.Range("B11").Value = 1757
$ws1.Range("B12").Value = 1147
$ws1.Range("B13").Value = 1776
$ws1.Range("B15").Value = 75978

# ---------------------------------------------------------------------------
# Sheet 2: "Cases by Gender"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = 25441
$ws2.Range("B3").Value = 49630
$ws2.Range("B4").Value = 907
$ws2.Range("B5").Value = 75978

# ---------------------------------------------------------------------------
# Sheet 3: "Cases by RaceEthnicity"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = 25714
$ws3.Range("B3").Value = 12715
$ws3.Range("B4").Value = 27770
$ws3.Range("B5").Value = 938
$ws3.Range("B6").Value = 438
$ws3.Range("B7").Value = 8403
$ws3.Range("B8").Value = 75978

# ---------------------------------------------------------------------------
# Sheet 4: "Fatalities by Age Group" (totals unchanged, formula -> literal)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B15").Value = 44134

# ---------------------------------------------------------------------------
# Sheet 5: "Fatalities by Gender" (totals unchanged, formula -> literal)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B5").Value = 44134

# ---------------------------------------------------------------------------
# Sheet 6: "Fatalities by Race-Ethnicity" (no data changes)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# ---------------------------------------------------------------------------
# Update selections on each sheet to match the saved view state, and make
# sheet 6 ("Fatalities by Race-Ethnicity") the active/selected tab.
# ---------------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("B2:B15").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B2:B5").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("B2:B8").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("B2:B15").Select() | Out-Null

$ws5.Activate() | Out-Null
$ws5.Range("B2:B5").Select() | Out-Null

$ws6.Activate() | Out-Null
$ws6.Range("C13").Select() | Out-Null
